$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 93, shifting existing rows 93-96 down to 94-97.
$ws.Rows.Item(93).Insert()

# Populate the new row 93 with the new weekly record.
$ws.Range("A93").Value = 8
$ws.Range("B93").Value = "Terminal La Palmera de La Serena"
$ws.Range("C93").Value = "Coquimbo"
$ws.Range("D93").Value = 44509
$ws.Range("E93").Value = 4
$ws.Range("F93").Value = 100112001
$ws.Range("G93").Value = "Berenjena"
$ws.Range("H93").Value = "Sin especificar"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 400
$ws.Range("K93").Value = 8000
$ws.Range("L93").Value = 9000
$ws.Range("M93").Value = 8500
$ws.Range("N93").Value = "$/caja 60 unidades"
$ws.Range("O93").Value = "Región de Arica y Parinacota"
$ws.Range("P93").Value = 142
$ws.Range("Q93").Value = 60
$ws.Range("R93").Value = "Hortaliza"
